# Auto-generated: apply cryptocurrency price/volume refresh per commit
# "Updated cryptos list on Fri Mar 24 19:29:26 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.820.95'
$ws.Range("E2").Value = '  -0.61%  '
$ws.Range("D3").Value = '1.762.65'
$ws.Range("E3").Value = '  -2.50%  '
$ws.Range("E4").Value = '  -0.58%  '
$ws.Range("D5").Value = '321.64'
$ws.Range("E5").Value = '  -2.42%  '
$ws.Range("D6").Value = '0.9986'
$ws.Range("E6").Value = '  -0.43%  '
$ws.Range("D7").Value = '0.4248'
$ws.Range("E7").Value = '  -3.91%  '
$ws.Range("D8").Value = '0.3634'
$ws.Range("E8").Value = '  -2.37%  '
$ws.Range("D9").Value = '42.48'
$ws.Range("E9").Value = '  -5.05%  '
$ws.Range("D10").Value = '0.07483'
$ws.Range("E10").Value = '  -2.88%  '
$ws.Range("D11").Value = '1.088'
$ws.Range("E11").Value = '  -2.66%  '
$ws.Range("D12").Value = '0.9989'
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("D13").Value = '20.68'
$ws.Range("E13").Value = '  -5.56%  '
$ws.Range("D14").Value = '6.073'
$ws.Range("E14").Value = '  -3.42%  '
$ws.Range("D15").Value = '7.291'
$ws.Range("E15").Value = '  -2.37%  '
$ws.Range("D16").Value = '1.767.73'
$ws.Range("E16").Value = '  -2.67%  '
$ws.Range("D17").Value = '91.02'
$ws.Range("E17").Value = '  -2.59%  '
$ws.Range("E18").Value = '  -2.53%  '
$ws.Range("D19").Value = '0.06370'
$ws.Range("E19").Value = '  -1.63%  '
$ws.Range("D20").Value = '0.9990'
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").Value = '17.04'
$ws.Range("E21").Value = '  -2.33%  '
$ws.Range("D22").Value = '5.931'
$ws.Range("E22").Value = '  -5.16%  '
$ws.Range("D23").Value = '27.840.84'
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("D24").Value = '11.21'
$ws.Range("E24").Value = '  -3.90%  '
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("D26").Value = '157.53'
$ws.Range("E26").Value = '  +1.26%  '
$ws.Range("D27").Value = '20.22'
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("D28").Value = '1.966.80'
$ws.Range("E28").Value = '  -2.90%  '
$ws.Range("D29").Value = '2.133'
$ws.Range("E29").Value = '  -8.04%  '
$ws.Range("D30").Value = '124.24'
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("D31").Value = '1.115'
$ws.Range("E31").Value = '  -6.97%  '
$ws.Range("E32").Value = '  +0.48%  '
$ws.Range("D33").Value = '5.558'
$ws.Range("E33").Value = '  -4.81%  '
$ws.Range("D34").Value = '0.08859'
$ws.Range("E34").Value = '  -4.07%  '
$ws.Range("D35").Value = '12.22'
$ws.Range("E35").Value = '  -6.41%  '
$ws.Range("D36").Value = '0.02292'
$ws.Range("E36").Value = '  -2.03%  '
$ws.Range("D37").Value = '0.2102'
$ws.Range("E37").Value = '  -3.07%  '
$ws.Range("D38").Value = '0.06047'
$ws.Range("E38").Value = '  -2.48%  '
$ws.Range("E39").Value = '  -3.80%  '
$ws.Range("D40").Value = '0.6321'
$ws.Range("E40").Value = '  -3.63%  '
$ws.Range("D41").Value = '1.174'
$ws.Range("E41").Value = '  -1.69%  '
$ws.Range("D42").Value = '0.9979'
$ws.Range("E42").Value = '  -0.49%  '
$ws.Range("D43").Value = '7.868'
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("D44").Value = '1.394'
$ws.Range("E44").Value = '  +0.84%  '
$ws.Range("D45").Value = '13.36'
$ws.Range("E45").Value = '  -4.48%  '
$ws.Range("D46").Value = '0.5868'
$ws.Range("E46").Value = '  -3.33%  '
$ws.Range("D47").Value = '3.682'
$ws.Range("E47").Value = '  -2.18%  '
$ws.Range("D48").Value = '1.984'
$ws.Range("E48").Value = '  -2.34%  '
$ws.Range("D49").Value = '123.05'
$ws.Range("E49").Value = '  -2.87%  '
$ws.Range("D50").Value = '1.184'
$ws.Range("E50").Value = '  +2.89%  '
$ws.Range("D51").Value = '0.06831'
$ws.Range("E51").Value = '  -2.05%  '
